$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added at the top of the table (weekly log).
# Insert a fresh row above the current row 2, pushing all existing
# records down by one, then clear the formatting Excel auto-copies from
# the header row so the new row starts out unstyled like the rest of the
# data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Cells.Item(2, 1).Value  = 1
$ws.Cells.Item(2, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value  = 44616
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value  = 15
$ws.Cells.Item(2, 6).Value  = 100112044
$ws.Cells.Item(2, 7).Value  = "Perejil"
$ws.Cells.Item(2, 8).Value  = "Sin especificar"
$ws.Cells.Item(2, 9).Value  = "Primera"
$ws.Cells.Item(2, 10).Value = 270
$ws.Cells.Item(2, 11).Value = 1300
$ws.Cells.Item(2, 12).Value = 1500
$ws.Cells.Item(2, 13).Value = 1400
$ws.Cells.Item(2, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 700
$ws.Cells.Item(2, 17).Value = 2
$ws.Cells.Item(2, 18).Value = "Hortaliza"
